$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column price cells whose new values look numeric,
# so Excel stores them as text (matching the source data) instead of coercing to a Double.
$textCells = @("D5","D6","D12","D13","D19","D22","D24","D29","D30","D32","D34","D35","D37","D38","D40","D41","D42","D46","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '68.176.04'
$ws.Cells.Item(2, 5).Value = '  -0.63%  '
$ws.Cells.Item(3, 4).Value = '2.641.90'
$ws.Cells.Item(3, 5).Value = '  -0.34%  '
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).Value = '597.41'
$ws.Cells.Item(6, 4).Value = '156.37'
$ws.Cells.Item(6, 5).Value = '  +1.09%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 5).Value = '  -0.78%  '
$ws.Cells.Item(9, 5).Value = '  +2.39%  '
$ws.Cells.Item(11, 5).Value = '  +0.52%  '
$ws.Cells.Item(12, 4).Value = '0.350'
$ws.Cells.Item(12, 5).Value = '  +0.46%  '
$ws.Cells.Item(13, 4).Value = '27.92'
$ws.Cells.Item(13, 5).Value = '  +0.06%  '
$ws.Cells.Item(14, 5).Value = '  +0.52%  '
$ws.Cells.Item(15, 4).Value = '3.122.61'
$ws.Cells.Item(15, 5).Value = '  -0.34%  '
$ws.Cells.Item(16, 4).Value = '68.293.76'
$ws.Cells.Item(16, 5).Value = '  -0.37%  '
$ws.Cells.Item(17, 4).Value = '2.651.08'
$ws.Cells.Item(17, 5).Value = '  -0.17%  '
$ws.Cells.Item(19, 4).Value = '362.56'
$ws.Cells.Item(19, 5).Value = '  -1.98%  '
$ws.Cells.Item(20, 5).Value = '  -1.61%  '
$ws.Cells.Item(21, 5).Value = '  +3.13%  '
$ws.Cells.Item(22, 4).Value = '4.78'
$ws.Cells.Item(22, 5).Value = '  -1.68%  '
$ws.Cells.Item(23, 5).Value = '  -3.05%  '
$ws.Cells.Item(24, 4).Value = '75.30'
$ws.Cells.Item(24, 5).Value = '  +3.55%  '
$ws.Cells.Item(25, 5).Value = '  +0.00%  '
$ws.Cells.Item(26, 5).Value = '  -2.46%  '
$ws.Cells.Item(27, 5).Value = '  -0.03%  '
$ws.Cells.Item(28, 5).Value = '  -1.16%  '
$ws.Cells.Item(29, 4).Value = '0.997'
$ws.Cells.Item(29, 5).Value = '  -0.17%  '
$ws.Cells.Item(30, 4).Value = '555.20'
$ws.Cells.Item(30, 5).Value = '  -3.29%  '
$ws.Cells.Item(31, 5).Value = '  +0.47%  '
$ws.Cells.Item(32, 4).Value = '1.41'
$ws.Cells.Item(32, 5).Value = '  -0.96%  '
$ws.Cells.Item(33, 5).Value = '  +0.12%  '
$ws.Cells.Item(34, 4).Value = '0.999'
$ws.Cells.Item(34, 5).Value = '  -0.02%  '
$ws.Cells.Item(35, 4).Value = '0.127'
$ws.Cells.Item(35, 5).Value = '  -2.12%  '
$ws.Cells.Item(36, 5).Value = '  -0.04%  '
$ws.Cells.Item(37, 4).Value = '160.79'
$ws.Cells.Item(37, 5).Value = '  +1.57%  '
$ws.Cells.Item(38, 4).Value = '19.61'
$ws.Cells.Item(38, 5).Value = '  +1.88%  '
$ws.Cells.Item(39, 5).Value = '  +0.91%  '
$ws.Cells.Item(40, 4).Value = '1.87'
$ws.Cells.Item(40, 5).Value = '  -3.44%  '
$ws.Cells.Item(41, 4).Value = '5.31'
$ws.Cells.Item(41, 5).Value = '  -1.63%  '
$ws.Cells.Item(42, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(42, 4).Value = '17.79'
$ws.Cells.Item(42, 5).Value = '  +0.30%  '
$ws.Cells.Item(43, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(43, 4).Value = '0.0₆0333'
$ws.Cells.Item(43, 5).Value = '  +3.42%  '
$ws.Cells.Item(44, 5).Value = '  -2.37%  '
$ws.Cells.Item(45, 5).Value = '  +0.03%  '
$ws.Cells.Item(46, 4).Value = '158.48'
$ws.Cells.Item(47, 5).Value = '  -0.58%  '
$ws.Cells.Item(48, 4).Value = '21.91'
$ws.Cells.Item(48, 5).Value = '  -0.19%  '
$ws.Cells.Item(49, 4).Value = '0.0781'
$ws.Cells.Item(49, 5).Value = '  +0.22%  '
$ws.Cells.Item(50, 5).Value = '  -2.22%  '
$ws.Cells.Item(51, 5).Value = '  -0.69%  '
